$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New account-statement rows (periods 1809-1903) for the three workers.
# Row layout: C=Doc, D=Name, E=Period, F=ValorMora (G=SalarioBasico unchanged)
$rows = @(
    @{ r = 16; c = "73550754"; d = "ANGEL ESTEBAN MARTINEZ CARDENAS"; e = "1903"; f = 26041 },
    @{ r = 17; c = "73550754"; d = "ANGEL ESTEBAN MARTINEZ CARDENAS"; e = "1902"; f = 31249 },
    @{ r = 18; c = "73550754"; d = "ANGEL ESTEBAN MARTINEZ CARDENAS"; e = "1901"; f = 31249 },
    @{ r = 19; c = "73550754"; d = "ANGEL ESTEBAN MARTINEZ CARDENAS"; e = "1812"; f = 31249 },
    @{ r = 20; c = "73550754"; d = "ANGEL ESTEBAN MARTINEZ CARDENAS"; e = "1811"; f = 31249 },
    @{ r = 21; c = "73550754"; d = "ANGEL ESTEBAN MARTINEZ CARDENAS"; e = "1810"; f = 31249 },
    @{ r = 22; c = "73550754"; d = "ANGEL ESTEBAN MARTINEZ CARDENAS"; e = "1809"; f = 31249 },
    @{ r = 23; c = "73546271"; d = "DAVID FRANCISCO ROMERO JARABA";   e = "1903"; f = 26041 },
    @{ r = 24; c = "73546271"; d = "DAVID FRANCISCO ROMERO JARABA";   e = "1902"; f = 31249 },
    @{ r = 25; c = "73546271"; d = "DAVID FRANCISCO ROMERO JARABA";   e = "1901"; f = 31249 },
    @{ r = 26; c = "73546271"; d = "DAVID FRANCISCO ROMERO JARABA";   e = "1812"; f = 31249 },
    @{ r = 27; c = "73546271"; d = "DAVID FRANCISCO ROMERO JARABA";   e = "1811"; f = 31249 },
    @{ r = 28; c = "73546271"; d = "DAVID FRANCISCO ROMERO JARABA";   e = "1810"; f = 31249 },
    @{ r = 29; c = "73429088"; d = "ISAITH ENRIQUE BOBADILLA ARDILA"; e = "1903"; f = 26041 },
    @{ r = 30; c = "73429088"; d = "ISAITH ENRIQUE BOBADILLA ARDILA"; e = "1902"; f = 31249 },
    @{ r = 31; c = "73429088"; d = "ISAITH ENRIQUE BOBADILLA ARDILA"; e = "1901"; f = 31249 },
    @{ r = 32; c = "73429088"; d = "ISAITH ENRIQUE BOBADILLA ARDILA"; e = "1812"; f = 31249 },
    @{ r = 33; c = "73429088"; d = "ISAITH ENRIQUE BOBADILLA ARDILA"; e = "1811"; f = 31249 },
    @{ r = 34; c = "73429088"; d = "ISAITH ENRIQUE BOBADILLA ARDILA"; e = "1810"; f = 31249 }
)

foreach ($row in $rows) {
    $ws.Cells.Item($row.r, 3).Value = $row.c
    $ws.Cells.Item($row.r, 4).Value = $row.d
    $ws.Cells.Item($row.r, 5).Value = $row.e
    $ws.Cells.Item($row.r, 6).Value = $row.f
}

$ws.Columns.Item(2).ColumnWidth = 18.54296875
$ws.Columns.Item(3).ColumnWidth = 16.7265625
$ws.Columns.Item(5).ColumnWidth = 13.54296875
$ws.Columns.Item(6).ColumnWidth = 10.1796875
$ws.Columns.Item(7).ColumnWidth = 14.36328125
$ws.Columns.Item(8).ColumnWidth = 19.36328125
$ws.Columns.Item(9).ColumnWidth = 18.08984375
$ws.Columns.Item(10).ColumnWidth = 15
